# Apply the "Add file python for process audio" edit to the plan workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: col A gets a new custom width, col C gets a wider custom
# width (values chosen so the engine's pixel-quantized ColumnWidth setter
# lands as close as possible to the target 18.75 / 19.5 character widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 18.714285714285715

# ---------------------------------------------------------------------------
# C2 already has vertical-center alignment (style 2); adding wrap text on
# top of that produces a brand-new cellXf (vertical=center + wrapText).
# ---------------------------------------------------------------------------
$ws.Range("C2").WrapText = $true

# C3:C6 get plain wrap text, which reuses the existing wrap-only style (1).
$ws.Range("C3").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("C6").WrapText = $true

# ---------------------------------------------------------------------------
# New rows 7 and 8 describing the "gape" screens.
# The order in which brand-new strings are first written controls their
# position in sharedStrings.xml, so we deliberately write column A for both
# rows before column C, and row 8's text before row 7's text.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "./download/gape"
$ws.Range("B7").Value = "get"
$ws.Range("A8").Value = "./script/gape"
$ws.Range("B8").Value = "get"
$ws.Range("C8").Value = "load màn hình để upload thêm script"
$ws.Range("C7").Value = "load màn hình để download dữ liệu"

$ws.Range("C7").WrapText = $true
$ws.Range("C8").WrapText = $true

$ws.Rows.Item(7).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 31.5

# ---------------------------------------------------------------------------
# Update the active selection to D6 (was D2).
# ---------------------------------------------------------------------------
$ws.Range("D6").Select() | Out-Null
